# Update cryptocurrency price/volume figures in cryptos.xlsx
# (mirrors "Updated cryptos list ... with GitHub Actions" data refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.533.61'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '3.388.54'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.474'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("D12").Value = '3.969.05'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D15").Value = '3.404.63'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '61.553.08'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000112'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.194'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.83%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.93'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.23%  '
$ws.Range("D37").Value = '3.422.43'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0769'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.30%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("D45").Value = '2.454.78'
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.26%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.41%  '
$ws.Range("E51").Value = '  -1.60%  '
